# Rainfall workbook update: add 2022 monthly data, fix the August/December
# 2021 figures ("unscale" them), highlight the two corrected 2021 month cells
# in red, and leave the selection on the new last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / "unscale" the existing August 2021 and December 2021 rows ---------

# Row 237: August 2021
$ws.Range("C237").Value = 4
$ws.Range("D237").Value = 4.7

# Row 241: December 2021
$ws.Range("C241").Value = 0.3
$ws.Range("D241").Value = 0.3

# Highlight the corrected month cells (column A) in red to flag the change.
$ws.Range("A237").Interior.Color = 255
$ws.Range("A241").Interior.Color = 255

# --- Append the 2022 monthly rows (rows 243-253) -----------------------------

$months2022 = @(
    @(243, "Feb",  0.1, 0.2),
    @(244, "Mar",  0,   0),
    @(245, "Apr",  0.7, 0.3),
    @(246, "May",  7.8, 7.5),
    @(247, "Jun",  0.2, 1.4),
    @(248, "Jul",  4.1, 1.6),
    @(249, "Aug",  6,   4),
    @(250, "Sept", 1,   0.7),
    @(251, "Oct",  2.7, 3.8),
    @(252, "Nov",  0.7, 0.8),
    @(253, "Dec",  1.7, 1.4)
)

foreach ($row in $months2022) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = 2022
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

# --- Move the on-screen selection down to the new last row of data ----------

$ws.Range("F252").Select() | Out-Null
